$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 624
$ws1.Range("F3").Value = 614
$ws1.Range("F4").Value = 898
$ws1.Range("F5").Value = 653
$ws1.Range("F6").Value = 799
$ws1.Range("F7").Value = 369
$ws1.Range("F8").Value = 572
$ws1.Range("F10").Value = 1156
$ws1.Range("F11").Value = 590
$ws1.Range("F13").Value = 470
$ws1.Range("F15").Value = 307
$ws1.Range("F17").Value = 71
$ws1.Range("F19").Value = 34
$ws1.Range("F20").Value = 536
$ws1.Range("F22").Value = 524

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 56
$ws2.Range("F8").Value = 172
$ws2.Range("F9").Value = 201
$ws2.Range("F13").Value = 29

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 624
$ws4.Range("F5").Value = 56
$ws4.Range("F7").Value = 614
$ws4.Range("F8").Value = 898
$ws4.Range("F9").Value = 653
$ws4.Range("F10").Value = 799
$ws4.Range("F11").Value = 369
$ws4.Range("F12").Value = 572
$ws4.Range("F14").Value = 1156
$ws4.Range("F15").Value = 590
$ws4.Range("F19").Value = 470
$ws4.Range("F22").Value = 172
$ws4.Range("F23").Value = 307
$ws4.Range("F25").Value = 71
$ws4.Range("F26").Value = 201
$ws4.Range("F31").Value = 29
$ws4.Range("F32").Value = 34
$ws4.Range("F33").Value = 536
$ws4.Range("F35").Value = 524
